$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 34.408928
$ws.Range("H2").Value = 68.81785599999999
$ws.Range("I2").Value = 0.02296116112547488
$ws.Range("J2").Value = 0.01585611315973826
$ws.Range("M2").Value = 7.369448
$ws.Range("N2").Value = 14.738896
$ws.Range("O2").Value = 0.7452608427984224
$ws.Range("P2").Value = 0.661061693471796
$ws.Range("Q2").Value = 253.574805631744
$ws.Range("R2").Value = 1014.299222526976
$ws.Range("S2").Value = 0.01711205429200179
$ws.Range("T2").Value = 0.01048186901725701
$ws.Range("G3").Value = 34.408928
$ws.Range("H3").Value = 68.81785599999999
$ws.Range("I3").Value = 0.02296116112547488
$ws.Range("J3").Value = 0.01585611315973826
$ws.Range("O3").Value = 0.01116592909756377
$ws.Range("P3").Value = 0.01485661309677453
$ws.Range("Q3").Value = 3.799204436906666
$ws.Range("R3").Value = 22.79522662143999
$ws.Range("S3").Value = 0.0002563826971247901
$ws.Range("T3").Value = 0.0002355681384329064
$ws.Range("G4").Value = 34.408928
$ws.Range("H4").Value = 68.81785599999999
$ws.Range("I4").Value = 0.02296116112547488
$ws.Range("J4").Value = 0.01585611315973826
$ws.Range("M4").Value = 0.084843
$ws.Range("N4").Value = 0.254529
$ws.Range("O4").Value = 0.008580040959044227
$ws.Range("P4").Value = 0.0114160091622658
$ws.Range("Q4").Value = 2.919356678304
$ws.Range("R4").Value = 17.516140069824
$ws.Range("S4").Value = 0.0001970077029237885
$ws.Range("T4").Value = 0.0001810135331094953
$ws.Range("G5").Value = 34.408928
$ws.Range("H5").Value = 68.81785599999999
$ws.Range("I5").Value = 0.02296116112547488
$ws.Range("J5").Value = 0.01585611315973826
$ws.Range("M5").Value = 2.32371
$ws.Range("N5").Value = 6.97113
$ws.Range("O5").Value = 0.2349931871449696
$ws.Range("P5").Value = 0.3126656842691638
$ws.Range("Q5").Value = 79.95637008288
$ws.Range("R5").Value = 479.73822049728
$ws.Range("S5").Value = 0.005395716433424521
$ws.Range("T5").Value = 0.004957662470938857
$ws.Range("I6").Value = 0.07707058947984194
$ws.Range("J6").Value = 0.07983306994376788
$ws.Range("M6").Value = 7.369448
$ws.Range("N6").Value = 14.738896
$ws.Range("O6").Value = 0.7452608427984224
$ws.Range("P6").Value = 0.661061693471796
$ws.Range("Q6").Value = 851.1398722598641
$ws.Range("R6").Value = 5106.839233559184
$ws.Range("S6").Value = 0.05743769247071823
$ws.Range("T6").Value = 0.05277458441207953
$ws.Range("I7").Value = 0.07707058947984194
$ws.Range("J7").Value = 0.07983306994376788
$ws.Range("O7").Value = 0.01116592909756377
$ws.Range("P7").Value = 0.01485661309677453
$ws.Range("S7").Value = 0.0008605647376393594
$ws.Range("T7").Value = 0.001186049032482299
$ws.Range("I8").Value = 0.07707058947984194
$ws.Range("J8").Value = 0.07983306994376788
$ws.Range("M8").Value = 0.084843
$ws.Range("N8").Value = 0.254529
$ws.Range("O8").Value = 0.008580040959044227
$ws.Range("P8").Value = 0.0114160091622658
$ws.Range("Q8").Value = 9.799005323349
$ws.Range("R8").Value = 88.19104791014101
$ws.Range("S8").Value = 0.0006612688144747269
$ws.Range("T8").Value = 0.0009113750579298606
$ws.Range("I9").Value = 0.07707058947984194
$ws.Range("J9").Value = 0.07983306994376788
$ws.Range("M9").Value = 2.32371
$ws.Range("N9").Value = 6.97113
$ws.Range("O9").Value = 0.2349931871449696
$ws.Range("P9").Value = 0.3126656842691638
$ws.Range("Q9").Value = 268.3786129665301
$ws.Range("R9").Value = 2415.40751669877
$ws.Range("S9").Value = 0.01811106345700962
$ws.Range("T9").Value = 0.0249610614412762
$ws.Range("G10").Value = 366.2779236666667
$ws.Range("H10").Value = 1098.833771
$ws.Range("I10").Value = 0.2444181469999509
$ws.Range("J10").Value = 0.253178951357013
$ws.Range("M10").Value = 7.369448
$ws.Range("N10").Value = 14.738896
$ws.Range("O10").Value = 0.7452608427984224
$ws.Range("P10").Value = 0.661061693471796
$ws.Range("Q10").Value = 2699.26611200947
$ws.Range("R10").Value = 16195.59667205682
$ws.Range("S10").Value = 0.1821552742284121
$ws.Range("T10").Value = 0.1673669063354805
$ws.Range("G11").Value = 366.2779236666667
$ws.Range("H11").Value = 1098.833771
$ws.Range("I11").Value = 0.2444181469999509
$ws.Range("J11").Value = 0.253178951357013
$ws.Range("O11").Value = 0.01116592909756377
$ws.Range("P11").Value = 0.01485661309677453
$ws.Range("Q11").Value = 40.44196647844889
$ws.Range("R11").Value = 363.97769830604
$ws.Range("S11").Value = 0.002729155699559371
$ws.Range("T11").Value = 0.00376138172455824
$ws.Range("G12").Value = 366.2779236666667
$ws.Range("H12").Value = 1098.833771
$ws.Range("I12").Value = 0.2444181469999509
$ws.Range("J12").Value = 0.253178951357013
$ws.Range("M12").Value = 0.084843
$ws.Range("N12").Value = 0.254529
$ws.Range("O12").Value = 0.008580040959044227
$ws.Range("P12").Value = 0.0114160091622658
$ws.Range("Q12").Value = 31.076117877651
$ws.Range("R12").Value = 279.685060898859
$ws.Range("S12").Value = 0.002097117712393271
$ws.Range("T12").Value = 0.002890293228384507
$ws.Range("G13").Value = 366.2779236666667
$ws.Range("H13").Value = 1098.833771
$ws.Range("I13").Value = 0.2444181469999509
$ws.Range("J13").Value = 0.253178951357013
$ws.Range("M13").Value = 2.32371
$ws.Range("N13").Value = 6.97113
$ws.Range("O13").Value = 0.2349931871449696
$ws.Range("P13").Value = 0.3126656842691638
$ws.Range("Q13").Value = 851.1236740034701
$ws.Range("R13").Value = 7660.113066031231
$ws.Range("S13").Value = 0.05743659935958616
$ws.Range("T13").Value = 0.07916037006858979
$ws.Range("G14").Value = 121.157162
$ws.Range("H14").Value = 242.314324
$ws.Range("I14").Value = 0.08084846811232432
$ws.Range("J14").Value = 0.05583090733848903
$ws.Range("M14").Value = 7.369448
$ws.Range("N14").Value = 14.738896
$ws.Range("O14").Value = 0.7452608427984224
$ws.Range("P14").Value = 0.661061693471796
$ws.Range("Q14").Value = 892.861405186576
$ws.Range("R14").Value = 3571.445620746304
$ws.Range("S14").Value = 0.06025319748435221
$ws.Range("T14").Value = 0.03690767415324848
$ws.Range("G15").Value = 121.157162
$ws.Range("H15").Value = 242.314324
$ws.Range("I15").Value = 0.08084846811232432
$ws.Range("J15").Value = 0.05583090733848903
$ws.Range("O15").Value = 0.01116592909756377
$ws.Range("P15").Value = 0.01485661309677453
$ws.Range("Q15").Value = 13.37736611362667
$ws.Range("R15").Value = 80.26419668176
$ws.Range("S15").Value = 0.0009027482625888588
$ws.Range("T15").Value = 0.0008294581891698013
$ws.Range("G16").Value = 121.157162
$ws.Range("H16").Value = 242.314324
$ws.Range("I16").Value = 0.08084846811232432
$ws.Range("J16").Value = 0.05583090733848903
$ws.Range("M16").Value = 0.084843
$ws.Range("N16").Value = 0.254529
$ws.Range("O16").Value = 0.008580040959044227
$ws.Range("P16").Value = 0.0114160091622658
$ws.Range("Q16").Value = 10.279337095566
$ws.Range("R16").Value = 61.676022573396
$ws.Range("S16").Value = 0.0006936831678797237
$ws.Range("T16").Value = 0.0006373661497138037
$ws.Range("G17").Value = 121.157162
$ws.Range("H17").Value = 242.314324
$ws.Range("I17").Value = 0.08084846811232432
$ws.Range("J17").Value = 0.05583090733848903
$ws.Range("M17").Value = 2.32371
$ws.Range("N17").Value = 6.97113
$ws.Range("O17").Value = 0.2349931871449696
$ws.Range("P17").Value = 0.3126656842691638
$ws.Range("Q17").Value = 281.53410891102
$ws.Range("R17").Value = 1689.20465346612
$ws.Range("S17").Value = 0.01899883919750354
$ws.Range("T17").Value = 0.01745640884635695
$ws.Range("G18").Value = 672.2915446666666
$ws.Range("H18").Value = 2016.874634
$ws.Range("I18").Value = 0.4486217786379665
$ws.Range("J18").Value = 0.4647019579585521
$ws.Range("M18").Value = 7.369448
$ws.Range("N18").Value = 14.738896
$ws.Range("O18").Value = 0.7452608427984224
$ws.Range("P18").Value = 0.661061693471796
$ws.Range("Q18").Value = 4954.417579260677
$ws.Range("R18").Value = 29726.50547556406
$ws.Range("S18").Value = 0.3343402448454583
$ws.Range("T18").Value = 0.3071966632877398
$ws.Range("G19").Value = 672.2915446666666
$ws.Range("H19").Value = 2016.874634
$ws.Range("I19").Value = 0.4486217786379665
$ws.Range("J19").Value = 0.4647019579585521
$ws.Range("O19").Value = 0.01116592909756377
$ws.Range("P19").Value = 0.01485661309677453
$ws.Range("Q19").Value = 74.22995041846221
$ws.Range("R19").Value = 668.0695537661599
$ws.Range("S19").Value = 0.005009278971894484
$ws.Range("T19").Value = 0.00690389719470379
$ws.Range("G20").Value = 672.2915446666666
$ws.Range("H20").Value = 2016.874634
$ws.Range("I20").Value = 0.4486217786379665
$ws.Range("J20").Value = 0.4647019579585521
$ws.Range("M20").Value = 0.084843
$ws.Range("N20").Value = 0.254529
$ws.Range("O20").Value = 0.008580040959044227
$ws.Range("P20").Value = 0.0114160091622658
$ws.Range("Q20").Value = 57.039231524154
$ws.Range("R20").Value = 513.353083717386
$ws.Range("S20").Value = 0.003849193235833025
$ws.Range("T20").Value = 0.005305041809777688
$ws.Range("G21").Value = 672.2915446666666
$ws.Range("H21").Value = 2016.874634
$ws.Range("I21").Value = 0.4486217786379665
$ws.Range("J21").Value = 0.4647019579585521
$ws.Range("M21").Value = 2.32371
$ws.Range("N21").Value = 6.97113
$ws.Range("O21").Value = 0.2349931871449696
$ws.Range("P21").Value = 0.3126656842691638
$ws.Range("Q21").Value = 1562.21058525738
$ws.Range("R21").Value = 14059.89526731642
$ws.Range("S21").Value = 0.1054230615847808
$ws.Range("T21").Value = 0.1452963556663309
$ws.Range("G22").Value = 188.9396033333333
$ws.Range("H22").Value = 566.81881
$ws.Range("I22").Value = 0.1260798556444414
$ws.Range("J22").Value = 0.1305990002424397
$ws.Range("M22").Value = 7.369448
$ws.Range("N22").Value = 14.738896
$ws.Range("O22").Value = 0.7452608427984224
$ws.Range("P22").Value = 0.661061693471796
$ws.Range("Q22").Value = 1392.380581905627
$ws.Range("R22").Value = 8354.28349143376
$ws.Range("S22").Value = 0.09396237947747986
$ws.Range("T22").Value = 0.0863339962659907
$ws.Range("G23").Value = 188.9396033333333
$ws.Range("H23").Value = 566.81881
$ws.Range("I23").Value = 0.1260798556444414
$ws.Range("J23").Value = 0.1305990002424397
$ws.Range("O23").Value = 0.01116592909756377
$ws.Range("P23").Value = 0.01485661309677453
$ws.Range("Q23").Value = 20.86145140271111
$ws.Range("R23").Value = 187.7530626244
$ws.Range("S23").Value = 0.001407798728756909
$ws.Range("T23").Value = 0.00194025881742749
$ws.Range("G24").Value = 188.9396033333333
$ws.Range("H24").Value = 566.81881
$ws.Range("I24").Value = 0.1260798556444414
$ws.Range("J24").Value = 0.1305990002424397
$ws.Range("M24").Value = 0.084843
$ws.Range("N24").Value = 0.254529
$ws.Range("O24").Value = 0.008580040959044227
$ws.Range("P24").Value = 0.0114160091622658
$ws.Range("Q24").Value = 16.03020276561
$ws.Range("R24").Value = 144.27182489049
$ws.Range("S24").Value = 0.001081770325539691
$ws.Range("T24").Value = 0.001490919383350445
$ws.Range("G25").Value = 188.9396033333333
$ws.Range("H25").Value = 566.81881
$ws.Range("I25").Value = 0.1260798556444414
$ws.Range("J25").Value = 0.1305990002424397
$ws.Range("M25").Value = 2.32371
$ws.Range("N25").Value = 6.97113
$ws.Range("O25").Value = 0.2349931871449696
$ws.Range("P25").Value = 0.3126656842691638
$ws.Range("Q25").Value = 439.0408456617
$ws.Range("R25").Value = 3951.3676109553
$ws.Range("S25").Value = 0.02962790711266498
$ws.Range("T25").Value = 0.04083382577567111
